$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C27").Value = 2620
$ws.Range("C28").Value = 1294
$ws.Range("C29").Value = 755

$ws.Range("C33").Value = 2758
$ws.Range("C34").Value = 1347
$ws.Range("C35").Value = 783

$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Range("C34").Select()
